$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = 0.72782148605570085
$ws.Range("T1").Value = 0.91730158084496072
$ws.Range("W1").Value = 0.9889873971609594
$ws.Range("BA1").Value = 0.64442492166630028
$ws.Range("O3").Value = 0.96155529625095859
$ws.Range("AD3").Value = 0.94995509809755263
$ws.Range("AT3").Value = 0.80120285000941294
$ws.Range("BN3").Value = 0.88570554459801998
$ws.Range("AL4").Value = 0.65278123140829625
$ws.Range("M5").Value = 0.98589432589152182
$ws.Range("E6").Value = 0.61155666855954682
$ws.Range("P6").Value = 0.78871967281202249
$ws.Range("AK6").Value = 0.81109177047152747
$ws.Range("S7").Value = 0.75524472131840992
$ws.Range("AR7").Value = 0.84423922542799323
$ws.Range("BL7").Value = 0.69893676209304845
$ws.Range("BM9").Value = 0.69750950520373134
$ws.Range("BP9").Value = 0.71503501255272062
$ws.Range("M10").Value = 0.9387215695980724
$ws.Range("AF11").Value = 0.96019806770186134
$ws.Range("AN11").Value = 0.99894698459701048
$ws.Range("BA11").Value = 0.70803465626411288
$ws.Range("BO11").Value = 0.9884399892259419
$ws.Range("J12").Value = 0.74286651413008509
$ws.Range("Y12").Value = 0.6336619015535061
$ws.Range("AH12").Value = 0.75412747172302508
$ws.Range("V13").Value = 0.77665180725184679
$ws.Range("AF13").Value = 0.84305496518574574
$ws.Range("B14").Value = 0.66656740405269188
$ws.Range("BD14").Value = 0.87906529478895612
$ws.Range("BF14").Value = 0.9307096118971554
$ws.Range("BI14").Value = 0.65122610197715902
$ws.Range("V15").Value = 0.65856064904469303
$ws.Range("AH16").Value = 0.69382026019674115
$ws.Range("AW16").Value = 0.81977167738829837
$ws.Range("BE16").Value = 0.94071512483512221
$ws.Range("B17").Value = 0.95172784293012336
$ws.Range("I17").Value = 0.92459007331665677
$ws.Range("AF17").Value = 0.77862308543450209
$ws.Range("AQ17").Value = 0.98804721851213073
$ws.Range("R19").Value = 0.8343706790204648
$ws.Range("BB19").Value = 0.59042870435064021
$ws.Range("BF20").Value = 0.7471263489593406
$ws.Range("BI20").Value = 0.93581215215940228
$ws.Range("A21").Value = 0.90887050695020355
$ws.Range("AA22").Value = 0.82641247365060577
$ws.Range("U23").Value = 0.77337633661106109
$ws.Range("AG24").Value = 0.84895834325008757
$ws.Range("BI24").Value = 0.77520340945849409
$ws.Range("O25").Value = 0.73525145815112036
$ws.Range("X25").Value = 0.57688428615497989
$ws.Range("AA25").Value = 0.78514297290245505
$ws.Range("BN25").Value = 0.96651028912709913
$ws.Range("R26").Value = 0.6990534296227735
$ws.Range("T27").Value = 0.83243506691542546
$ws.Range("AQ27").Value = 0.838541492632495
$ws.Range("AS27").Value = 0.81060532329991764
$ws.Range("AA29").Value = 0.74971693109411297
$ws.Range("BE29").Value = 0.91557531435576811
$ws.Range("AS30").Value = 0.60350227671617862
$ws.Range("AG31").Value = 0.6619908933317652
$ws.Range("H32").Value = 0.74718467268915356
$ws.Range("P32").Value = 0.73571882319577564
$ws.Range("BI32").Value = 0.78330992715735215
$ws.Range("B33").Value = 0.70085407440568614
$ws.Range("AL33").Value = 0.94735018838405471
$ws.Range("J34").Value = 0.78489147129545445
$ws.Range("U35").Value = 0.79962128494926543
$ws.Range("AC35").Value = 0.91989305163533497
$ws.Range("AP35").Value = 0.97077266182799882
$ws.Range("K36").Value = 0.94314657479106634
$ws.Range("N36").Value = 0.73933914740844031
$ws.Range("Z36").Value = 0.86920057688244956
$ws.Range("I38").Value = 0.89378433002398094
$ws.Range("K38").Value = 0.87024647992725801
$ws.Range("AN38").Value = 0.83760810924844153
$ws.Range("G39").Value = 0.95991075825222505
$ws.Range("P39").Value = 0.95982966496184619
$ws.Range("AN39").Value = 0.55110446232193544
$ws.Range("AP40").Value = 0.95603424185853936
$ws.Range("AP41").Value = 0.57269317588996205
$ws.Range("M42").Value = 0.90851442015238004
$ws.Range("AE42").Value = 0.87348224865588175
$ws.Range("K43").Value = 0.86004611926222718
$ws.Range("E44").Value = 0.94744448537077641
$ws.Range("S44").Value = 0.83757080468764711
$ws.Range("AQ44").Value = 0.70169492304857584
$ws.Range("AJ45").Value = 0.76349743162165995
$ws.Range("J46").Value = 0.90191315487310653
$ws.Range("W46").Value = 0.88757670604094274
$ws.Range("AW47").Value = 0.89672757869870734
$ws.Range("H48").Value = 0.76125403519131618
$ws.Range("AU48").Value = 0.72868060886089747
$ws.Range("AB49").Value = 0.63980280351169339
$ws.Range("AK49").Value = 0.89866101532248188
$ws.Range("AZ49").Value = 0.95287847511856238
$ws.Range("BG49").Value = 0.9178376219986224
$ws.Range("R50").Value = 0.89943455239679249
$ws.Range("AI50").Value = 0.78108697129111881
$ws.Range("AV50").Value = 0.9960543591057569
$ws.Range("AB51").Value = 0.75492572709964256
$ws.Range("AR51").Value = 0.81876097669294312
$ws.Range("AX51").Value = 0.82409047338424224
$ws.Range("AZ51").Value = 0.61473782947470101
$ws.Range("BB51").Value = 0.86599206351005031
$ws.Range("BA54").Value = 0.77804103961797977
$ws.Range("R55").Value = 0.76631813269308791
$ws.Range("AW55").Value = 0.85391979141555874
$ws.Range("AB56").Value = 0.94221047001379166
$ws.Range("A57").Value = 0.79571703606231736
$ws.Range("C57").Value = 0.89291118955425586
$ws.Range("P58").Value = 0.63503467227820343
$ws.Range("BH59").Value = 0.7946272394951539
$ws.Range("D60").Value = 0.90542574973938339
$ws.Range("AS60").Value = 0.97834477105350159
$ws.Range("BI60").Value = 0.84893805236617759
$ws.Range("Q61").Value = 0.80660143830514985
$ws.Range("BK61").Value = 0.92152389976434801
$ws.Range("H62").Value = 0.60288088753892533
$ws.Range("AU62").Value = 0.81291255458088996
$ws.Range("AW62").Value = 0.87756479023246503
$ws.Range("AK63").Value = 0.96443754419350924
$ws.Range("BF63").Value = 0.87076352103681964
$ws.Range("AX64").Value = 0.79973656218311651
$ws.Range("BM64").Value = 0.81335356000266246
$ws.Range("M65").Value = 0.91465968169463718
$ws.Range("AB65").Value = 0.97277650762534784
$ws.Range("AO65").Value = 0.88992720922427715
$ws.Range("BO65").Value = 0.84965098595371591
$ws.Range("U66").Value = 0.70369828387440703
$ws.Range("BJ66").Value = 0.99621661241068549
$ws.Range("BL66").Value = 0.86829739644148585
$ws.Range("AZ67").Value = 0.94969865758218774
$ws.Range("BP67").Value = 0.98393164107267772
$ws.Range("B68").Value = 0.83097279652764688
$ws.Range("BL68").Value = 0.82224818544579426
